# Scheduled-runner update: refresh Universalis-sourced market-price and
# leve-profit figures (columns H-N) for a batch of leves across the
# ALC, BSM, GSM, LTW and WVR crafting sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 28: "The Writing Is Not on the Wall" (Enchanted Silver Ink)
$ws.Range("H28").Value = 1624444.9
$ws.Range("I28").Value = 2273903.5
$ws.Range("J28").Value = 798.375
$ws.Range("K28").Value = 2273903.5
$ws.Range("L28").Value = 798.375
$ws.Range("M28").Value = -2273418.5
$ws.Range("N28").Value = -1768.375

# Row 32: "Automata for the People" (Crab Oil)
$ws.Range("H32").Value = 648.0454999999999
$ws.Range("I32").Value = 487.75
$ws.Range("J32").Value = 683.6667
$ws.Range("K32").Value = 487.75
$ws.Range("L32").Value = 683.6667
$ws.Range("M32").Value = -161.75
$ws.Range("N32").Value = -1335.6667

# Row 33: "Glazed and Confused" (Clear Glass Lens)
$ws.Range("H33").Value = 5439
$ws.Range("I33").Value = 199
$ws.Range("K33").Value = 199
$ws.Range("M33").Value = 30

# Row 40: "Stuck in the Moment" (Horn Glue)
$ws.Range("H40").Value = 3000
$ws.Range("I40").Value = 4000
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 4000
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -3825
$ws.Range("N40").Value = -2350

# Row 53: "No Accounting for Waste" (Enchanted Electrum Ink)
$ws.Range("H53").Value = 422.13333
$ws.Range("I53").Value = 396.2
$ws.Range("J53").Value = 474
$ws.Range("K53").Value = 396.2
$ws.Range("L53").Value = 474
$ws.Range("M53").Value = 240.8
$ws.Range("N53").Value = -1748

# Row 76: "Warding Off Temptation" (Enchanted Hardsilver Ink)
$ws.Range("H76").Value = 6176906.5
$ws.Range("I76").Value = 6176906.5
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 6176906.5
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -6176591.5
$ws.Range("N76").ClearContents()

# Row 79: "The Garden of Arcane Delights (L)" (Enchanted Hardsilver Ink)
$ws.Range("H79").Value = 6176906.5
$ws.Range("I79").Value = 6176906.5
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 6176906.5
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -6175814.5
$ws.Range("N79").ClearContents()

# Row 98: "The Dotted Line" (Enchanted Durium Ink)
$ws.Range("H98").Value = 1590699.8
$ws.Range("I98").Value = 2060729.5
$ws.Range("J98").Value = 4349.625
$ws.Range("K98").Value = 2060729.5
$ws.Range("L98").Value = 4349.625
$ws.Range("M98").Value = -2059231.5
$ws.Range("N98").Value = -7345.625

# Row 107: "Another Man's Ink" (Enchanted Truegold Ink)
$ws.Range("H107").Value = 2764.3333
$ws.Range("I107").Value = 2474.25
$ws.Range("J107").Value = 2996.4
$ws.Range("K107").Value = 2474.25
$ws.Range("L107").Value = 2996.4
$ws.Range("M107").Value = -554.25
$ws.Range("N107").Value = -6836.4

# Row 113: "Amaro Kart" (Starch Glue)
$ws.Range("H113").Value = 2329.2307
$ws.Range("I113").Value = 2128
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2128
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1126
$ws.Range("N113").Value = -9508

# Row 116: "Growing Up" (Growth Formula Kappa)
$ws.Range("H116").Value = 3878.0667
$ws.Range("I116").Value = 8926.25
$ws.Range("J116").Value = 2042.3636
$ws.Range("K116").Value = 8926.25
$ws.Range("L116").Value = 2042.3636
$ws.Range("M116").Value = -5484.25
$ws.Range("N116").Value = -8926.363600000001

# Row 122: "Wishful Inking" (Enchanted High Durium Ink)
$ws.Range("H122").Value = 1590699.8
$ws.Range("I122").Value = 2060729.5
$ws.Range("J122").Value = 4349.625
$ws.Range("K122").Value = 6182188.5
$ws.Range("L122").Value = 13048.875
$ws.Range("M122").Value = -6179738.5
$ws.Range("N122").Value = -17948.875

# Row 132: "Fast-forwarding Flora" (Growth Formula Lambda)
$ws.Range("H132").Value = 12822246
$ws.Range("I132").Value = 2749126
$ws.Range("J132").Value = 83334090
$ws.Range("K132").Value = 8247378
$ws.Range("L132").Value = 250002270
$ws.Range("M132").Value = -8244848
$ws.Range("N132").Value = -250007330

# Row 137: "Cutting Edge of Culinary Quality" (Magnesia Whetstone)
$ws.Range("H137").Value = 1057.2979
$ws.Range("I137").Value = 1026.4348
$ws.Range("J137").Value = 1086.875
$ws.Range("K137").Value = 3079.3044
$ws.Range("L137").Value = 3260.625
$ws.Range("M137").Value = -529.3044
$ws.Range("N137").Value = -8360.625

# Row 138: "All-night Crafting" (Cunning Craftsman's Tisane)
$ws.Range("H138").Value = 2914.9792
$ws.Range("I138").Value = 1037.8846
$ws.Range("K138").Value = 3113.6538
$ws.Range("M138").Value = 2026.3462


$ws = $wb.Worksheets.Item("BSM")

# Row 105: "Ingot to Wing It" (Molybdenum Ingot)
$ws.Range("H105").Value = 4398.625
$ws.Range("I105").Value = 3698.3333
$ws.Range("J105").Value = 6499.5
$ws.Range("K105").Value = 3698.3333
$ws.Range("L105").Value = 6499.5
$ws.Range("M105").Value = -1951.3333
$ws.Range("N105").Value = -9993.5

# Row 107: "The Gold Experience" (Deepgold Nugget)
$ws.Range("H107").Value = 2721
$ws.Range("I107").Value = 2615.32
$ws.Range("J107").Value = 3381.5
$ws.Range("K107").Value = 2615.32
$ws.Range("L107").Value = 3381.5
$ws.Range("M107").Value = -695.3200000000002
$ws.Range("N107").Value = -7221.5


$ws = $wb.Worksheets.Item("GSM")

# Row 70: "Sky Is the Limit" (Mythrite Ingot)
$ws.Range("H70").Value = 21642732
$ws.Range("I70").Value = 40186044
$ws.Range("J70").Value = 8870.666999999999
$ws.Range("K70").Value = 40186044
$ws.Range("L70").Value = 8870.666999999999
$ws.Range("M70").Value = -40185774
$ws.Range("N70").Value = -9410.666999999999

# Row 73: "Hulls of Broken Dreams (L)" (Mythrite Ingot)
$ws.Range("H73").Value = 21642732
$ws.Range("I73").Value = 40186044
$ws.Range("J73").Value = 8870.666999999999
$ws.Range("K73").Value = 40186044
$ws.Range("L73").Value = 8870.666999999999
$ws.Range("M73").Value = -40185108
$ws.Range("N73").Value = -10742.667

# Row 113: "Copious Crystal Cannons" (Manasilver Nugget)
$ws.Range("H113").Value = 11905718
$ws.Range("I113").Value = 809.0909
$ws.Range("J113").Value = 25001118
$ws.Range("K113").Value = 809.0909
$ws.Range("L113").Value = 25001118
$ws.Range("M113").Value = 1360.9091
$ws.Range("N113").Value = -25005458

# Row 122: "Awarding Academic Excellence" (Ametrine)
$ws.Range("H122").Value = 2215.8
$ws.Range("I122").Value = 1923.7
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 5771.1
$ws.Range("L122").Value = 8400
$ws.Range("M122").Value = -3321.1
$ws.Range("N122").Value = -13300

# Row 126: "Gold Rush Order" (Phrygian Gold Ingot)
$ws.Range("H126").Value = 15154616
$ws.Range("I126").Value = 1678.6666
$ws.Range("J126").Value = 33338140
$ws.Range("K126").Value = 5035.9998
$ws.Range("L126").Value = 100014420
$ws.Range("M126").Value = -2565.9998
$ws.Range("N126").Value = -100019360


$ws = $wb.Worksheets.Item("LTW")

# Row 61: "Spelling Me Softly" (Raptor Leather)
$ws.Range("H61").Value = 1812.76
$ws.Range("I61").Value = 866.3570999999999
$ws.Range("J61").Value = 3017.2727
$ws.Range("K61").Value = 866.3570999999999
$ws.Range("L61").Value = 3017.2727
$ws.Range("M61").Value = -664.3570999999999
$ws.Range("N61").Value = -3421.2727

# Row 113: "Peace in Rest" (Atrociraptor Leather)
$ws.Range("H113").Value = 1812.76
$ws.Range("I113").Value = 866.3570999999999
$ws.Range("J113").Value = 3017.2727
$ws.Range("K113").Value = 866.3570999999999
$ws.Range("L113").Value = 3017.2727
$ws.Range("M113").Value = 1303.6429
$ws.Range("N113").Value = -7357.2727


$ws = $wb.Worksheets.Item("WVR")

# Row 81: "Where the Dragonflies, the Net Catches" (Crawler Silk)
$ws.Range("H81").Value = 12858.5
$ws.Range("J81").Value = 12858.5
$ws.Range("L81").Value = 25717
$ws.Range("N81").Value = -27839

# Row 84: "To Kill a Dragon on Nameday (L)" (Crawler Silk)
$ws.Range("H84").Value = 12858.5
$ws.Range("J84").Value = 12858.5
$ws.Range("L84").Value = 128585
$ws.Range("N84").Value = -139193
